$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario 2")

$ws.Range("C2").Value = 59.3
$ws.Range("D2").Value = 35.6
$ws.Range("E2").Value = 33.700000000000003
$ws.Range("F2").Value = 40.1

$ws.Activate()
$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollRow = 4
